# Updated 1st May 2024
# Applies the additions of new "Skills" rows (Spark/Databricks, Snowflake,
# Flink, Azure, AWS) to the Resources worksheet, widens columns D/E and
# updates the sheet view's selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -----------------------------------------------------
# Values are written in the same order the original author typed them so
# that the generated shared-string table matches the canonical one
# (new strings are appended to xl/sharedStrings.xml in assignment order).
$ws.Range("C43").Value = "Spark, Databricks"
$ws.Range("D44").Value = "Trendytech Ultimate Big Data Masters - Week 3, 4, 5, 6, 7, 8, 9, 10, 11"
$ws.Range("D51").Value = "Trendytech Ultimate Big Data Masters - Week 17, 18, 19, 20"
$ws.Range("D43").Value = "Trendytech Big Data Masters - Week 9, 10, 11, 12, 13, 14"
$ws.Range("D47").Value = "Trendytech Big Data Masters - Week 15, 16, 17"
$ws.Range("D48").Value = "Trendytech Ultimate Big Data Masters - Week 26, 27, 28, 29, 30"
$ws.Range("D49").Value = "Udemy - Prashant Pandey Stream Processing in Lakehouse course"
$ws.Range("D45").Value = "Scholarnest - Apace Spark Performance Tuning course"
$ws.Range("D52").Value = "Udemy - Prashant Pandey Master Databricks course"
$ws.Range("C55").Value = "Snowflake"
$ws.Range("C58").Value = "Flink"
$ws.Range("C61").Value = "Azure"
$ws.Range("C64").Value = "AWS"

# --- Column widths -------------------------------------------------------
# Column D grows from 57.140625 to 71.42578125 characters and a new
# column E is introduced at the old column-D width (35.7109375 chars).
$ws.Columns.Item(4).ColumnWidth = 70.6666666667
$ws.Columns.Item(5).ColumnWidth = 34.8333333333

# --- Sheet view ------------------------------------------------------------
# Scroll the view down and move the active selection to D54.
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D54").Select() | Out-Null
